$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.419.80"
$ws.Range("E2").Value = "  -0.05%  "

$ws.Range("D3").Value = "1.877.42"
$ws.Range("E3").Value = "  +0.16%  "

$ws.Range("D4").Value = "'0.9987"
$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").Value = "'0.7149"
$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("D6").Value = "'242.17"
$ws.Range("E6").Value = "  +0.43%  "

$ws.Range("D7").Value = "'0.9986"
$ws.Range("E7").Value = "  -0.28%  "

$ws.Range("D8").Value = "'0.3122"
$ws.Range("E8").Value = "  +1.02%  "

$ws.Range("D9").Value = "'0.07755"
$ws.Range("E9").Value = "  -1.75%  "

$ws.Range("D10").Value = "'24.99"
$ws.Range("E10").Value = "  -2.22%  "

$ws.Range("D11").Value = "'0.08383"
$ws.Range("E11").Value = "  +1.57%  "

$ws.Range("D12").Value = "1.902.28"
$ws.Range("E12").Value = "  +1.88%  "

$ws.Range("D13").Value = "'5.249"
$ws.Range("E13").Value = "  -0.21%  "

$ws.Range("D14").Value = "'0.7181"
$ws.Range("E14").Value = "  -1.10%  "

$ws.Range("D15").Value = "'91.54"
$ws.Range("E15").Value = "  +0.42%  "

$ws.Range("D16").Value = "29.405.59"
$ws.Range("E16").Value = "  -0.10%  "

$ws.Range("D17").Value = "'0.000008267"
$ws.Range("E17").Value = "  +5.58%  "

$ws.Range("D18").Value = "'5.984"
$ws.Range("E18").Value = "  +1.74%  "

$ws.Range("D19").Value = "'244.98"
$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").Value = "2.127.23"
$ws.Range("E20").Value = "  -0.13%  "

$ws.Range("D21").Value = "'13.22"
$ws.Range("E21").Value = "  -0.24%  "

$ws.Range("E22").Value = "  -0.22%  "

$ws.Range("D23").Value = "'7.941"
$ws.Range("E23").Value = "  -1.53%  "

$ws.Range("D24").Value = "'0.9988"
$ws.Range("E24").Value = "  -0.34%  "

$ws.Range("D25").Value = "'0.1634"
$ws.Range("E25").Value = "  +0.94%  "

$ws.Range("D26").Value = "'163.81"
$ws.Range("E26").Value = "  +0.82%  "

$ws.Range("D27").Value = "'9.032"
$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").Value = "'18.58"
$ws.Range("E28").Value = "  +1.55%  "

$ws.Range("E29").Value = "  +1.12%  "

$ws.Range("D30").Value = "'4.424"
$ws.Range("E30").Value = "  +0.44%  "

$ws.Range("D31").Value = "'1.301"
$ws.Range("E31").Value = "  -3.97%  "

$ws.Range("D32").Value = "'4.325"
$ws.Range("E32").Value = "  +5.48%  "

$ws.Range("D33").Value = "'0.05228"
$ws.Range("E33").Value = "  +0.58%  "

$ws.Range("D34").Value = "'1.929"
$ws.Range("E34").Value = "  -0.68%  "

$ws.Range("D35").Value = "'0.7714"
$ws.Range("E35").Value = "  +6.42%  "

$ws.Range("D36").Value = "'1.177"
$ws.Range("E36").Value = "  -1.48%  "

$ws.Range("D37").Value = "'2.684"
$ws.Range("E37").Value = "  +0.41%  "

$ws.Range("D38").Value = "'0.01869"
$ws.Range("E38").Value = "  +0.22%  "

$ws.Range("D39").Value = "'2.721"
$ws.Range("E39").Value = "  +0.97%  "

$ws.Range("D40").Value = "1.167.06"
$ws.Range("E40").Value = "  -1.70%  "

$ws.Range("D41").Value = "'6.397"
$ws.Range("E41").Value = "  +4.18%  "

$ws.Range("D42").Value = "'73.74"
$ws.Range("E42").Value = "  +1.30%  "

$ws.Range("D43").Value = "'0.8935"
$ws.Range("E43").Value = "  -1.46%  "

$ws.Range("D44").Value = "'104.01"
$ws.Range("E44").Value = "  +1.78%  "

$ws.Range("D45").Value = "'0.9986"
$ws.Range("E45").Value = "  -0.31%  "

$ws.Range("D46").Value = "2.023.44"
$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("D47").Value = "'1.805"
$ws.Range("E47").Value = "  +0.69%  "

$ws.Range("D48").Value = "'0.5195"
$ws.Range("E48").Value = "  -1.78%  "

$ws.Range("D49").Value = "'9.426"
$ws.Range("E49").Value = "  +1.30%  "

$ws.Range("D50").Value = "'0.4321"
$ws.Range("E50").Value = "  +0.26%  "

$ws.Range("D51").Value = "'7.074"
$ws.Range("E51").Value = "  +0.28%  "
